$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) onto the new
# header cells I1:J1, then overwrite their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Data rows: I column is always 1, J column mirrors the H column value.
for ($r = 2; $r -le 8; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
